# Add two new trade records (rows 5 & 6) to the CELG named-trade log.
# Columns: A=Principle, B=Start Principle, C=BuyPrice, D=SellPrice,
#          E=IsShortSell, F=Price Change %, G=Date (serial), H=Profitable
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - trade closed 2016-08-25 21:13:18, profitable
$ws.Cells.Item(5, 1).Value = 10146.58
$ws.Cells.Item(5, 2).Value = 9864.4599999999991
$ws.Cells.Item(5, 3).Value = 113.86
$ws.Cells.Item(5, 4).Value = 110.6
$ws.Cells.Item(5, 5).Value = $true
$ws.Cells.Item(5, 6).Value = -2.86
$ws.Cells.Item(5, 7).Value = 42607.884236111109
$ws.Cells.Item(5, 8).Value = $true

# Row 6 - trade closed 2016-08-26 14:47:30, profitable
$ws.Cells.Item(6, 1).Value = 10408.36
$ws.Cells.Item(6, 2).Value = 10146.58
$ws.Cells.Item(6, 3).Value = 110.77
$ws.Cells.Item(6, 4).Value = 107.91
$ws.Cells.Item(6, 5).Value = $true
$ws.Cells.Item(6, 6).Value = -2.58
$ws.Cells.Item(6, 7).Value = 42608.616319444445
$ws.Cells.Item(6, 8).Value = $true

# Column A's "best fit" width is recalculated by Excel once the wider
# values above are added (8.85546875 -> 9 characters). ColumnWidth uses
# a slightly different unit base than the stored <col width> value, so
# subtract the ~0.8333 offset to land exactly on the target width of 9.
$ws.Columns.Item(1).ColumnWidth = 8.166666666666666

